# "Generate Report for Archive"
# The localization status changed from "Ready for handoff" to "In Translation"
# for the tracked file. That status string is shared across the Overview
# sheet (columns E/F) and each per-locale sheet (column C). Updating the
# cell values rewrites the shared string used by every occurrence.
#
# Excel also narrows columns E/F (Overview) and C (zh-cn / de-de) to match
# the new, shorter status text - mirror that with ColumnWidth tweaks.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
# Closest ColumnWidth (character units) this runtime can resolve to the
# narrower post-edit column width used for the Status columns.
$newStatusColWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # Cast explicitly to string first: with boolean cell values (e.g.
            # "True"/"False" status flags) PowerShell's -eq coerces the
            # *string* operand to bool using the left-hand operand's type,
            # which would otherwise make "True" match any non-empty string.
            $cellText = [string]$cell.Value2
            if ($cellText -ceq $oldStatus) {
                $cell.Value = $newStatus
            }
        }
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = $newStatusColWidth
$overview.Range("F1").ColumnWidth = $newStatusColWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = $newStatusColWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = $newStatusColWidth
